$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update calculated values (rows 2-4, columns B-J) ---
$ws.Range("B2").Value = 0.018816
$ws.Range("C2").Value = 0.003857
$ws.Range("D2").Value = 3.070236
$ws.Range("E2").Value = 0.101042
$ws.Range("F2").Value = 0.6645
$ws.Range("G2").Value = 0.903
$ws.Range("H2").Value = 0.9990732
$ws.Range("I2").Value = 2.16066
$ws.Range("J2").Value = 65

$ws.Range("B3").Value = 0.03153
$ws.Range("C3").Value = 0.01031
$ws.Range("D3").Value = 2.89612
$ws.Range("E3").Value = 0.15028
$ws.Range("F3").Value = 0.9591
$ws.Range("G3").Value = 0.9694
$ws.Range("H3").Value = 0.9672969
$ws.Range("I3").Value = 2.524932
$ws.Range("J3").Value = 56

$ws.Range("B4").Value = 0.018127
$ws.Range("C4").Value = 0.005863
$ws.Range("D4").Value = 2.733211
$ws.Range("E4").Value = 0.196963
$ws.Range("F4").Value = 0.2053
$ws.Range("G4").Value = 0.7234
$ws.Range("H4").Value = 0.9992782
$ws.Range("I4").Value = 1.177094
$ws.Range("J4").Value = 124

# --- Clear leftover bold/font styling on I2:I4 (matches the source style) ---
$ws.Range("I2:I4").Font.Bold = $false

# --- Set explicit column widths to match the author's saved layout ---
$ws.Range("A1").EntireColumn.ColumnWidth = 22.666666666666668
$ws.Range("B1:E1").EntireColumn.ColumnWidth = 8.166666666666666
$ws.Range("F1:G1").EntireColumn.ColumnWidth = 6.166666666666667
$ws.Range("H1").EntireColumn.ColumnWidth = 12.333333333333334

# --- Restore the active selection cell ---
$ws.Range("H9").Select() | Out-Null
